$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.445647641019636;  C = 0.04103571897497393; D = 0.1496068669990043; E = 0.5333859586016987;  G = 2.169676185595313 }
    3 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    4 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    5 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    6 = @{ B = 0.1169995834814548; C = 0.3048912486333797;  D = 3.223369029078222;  E = 13.86384647080068;   G = 17.50910633199374 }
    7 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 13.86384647080068;   G = 19.48425592650926 }
    8 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
